# Auto-generated: applies the 88 cell-value edits from the commit diff
# (daily cryptos-list refresh) against Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.293.50'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '3.890.70'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''485.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = '''145.70'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '''0.621'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '''0.739'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').Value = '''0.176'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('D11').Value = '''0.0000351'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '''10.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').Value = '4.504.12'
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('D15').Value = '3.905.89'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').Value = '''14.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.45%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '''19.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '68.247.77'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '''429.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('D22').Value = '''3.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.28%  '
$ws.Range('D23').Value = '''14.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('D24').Value = '''12.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +18.22%  '
$ws.Range('D25').Value = '''88.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').Value = '''3.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.85%  '
$ws.Range('D27').Value = '''10.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.61%  '
$ws.Range('D28').Value = '''37.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.71%  '
$ws.Range('E29').Value = '  -3.26%  '
$ws.Range('D30').Value = '''720.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').Value = '''13.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.79%  '
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('E33').Value = '  +2.72%  '
$ws.Range('D34').Value = '''61.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.71%  '
$ws.Range('D35').Value = '''6.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.24%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = '''40.80'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0868'
$ws.Range('E37').Value = '  -4.08%  '
$ws.Range('D38').Value = '''0.400'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +18.27%  '
$ws.Range('E39').Value = '  -3.94%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = '''3.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.12%  '
$ws.Range('D42').Value = '''0.0494'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.91%  '
$ws.Range('E43').Value = '  +3.14%  '
$ws.Range('D44').Value = '''2.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('D47').Value = '''3.32'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.09%  '
$ws.Range('E48').Value = '  +31.08%  '
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('E50').Value = '  -2.61%  '
$ws.Range('D51').Value = '''144.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.49%  '
